$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1896
$ws.Range("I19").Value = 1839.6666
$ws.Range("K19").Value = 1839.6666
$ws.Range("M19").Value = -1664.6666
$ws.Range("H55").Value = 245.7
$ws.Range("J55").Value = 196.28572
$ws.Range("L55").Value = 196.28572
$ws.Range("N55").Value = -624.28572
$ws.Range("H64").Value = 7330.5835
$ws.Range("I64").Value = 5099.4
$ws.Range("J64").Value = 8924.286
$ws.Range("K64").Value = 5099.4
$ws.Range("L64").Value = 8924.286
$ws.Range("M64").Value = -4851.4
$ws.Range("N64").Value = -9420.286
$ws.Range("H67").Value = 7330.5835
$ws.Range("I67").Value = 5099.4
$ws.Range("J67").Value = 8924.286
$ws.Range("K67").Value = 5099.4
$ws.Range("L67").Value = 8924.286
$ws.Range("M67").Value = -4241.4
$ws.Range("N67").Value = -10640.286
$ws.Range("H69").Value = 28357.143
$ws.Range("I69").Value = 9083.166999999999
$ws.Range("K69").Value = 27249.501
$ws.Range("M69").Value = -26375.501
$ws.Range("H72").Value = 28357.143
$ws.Range("I72").Value = 9083.166999999999
$ws.Range("K72").Value = 81748.503
$ws.Range("M72").Value = -77380.503
$ws.Range("H80").Value = 553.5
$ws.Range("I80").Value = 606
$ws.Range("J80").Value = 474.75
$ws.Range("K80").Value = 1818
$ws.Range("L80").Value = 1424.25
$ws.Range("M80").Value = -820
$ws.Range("N80").Value = -3420.25
$ws.Range("H83").Value = 553.5
$ws.Range("I83").Value = 606
$ws.Range("J83").Value = 474.75
$ws.Range("K83").Value = 5454
$ws.Range("L83").Value = 4272.75
$ws.Range("M83").Value = -462
$ws.Range("N83").Value = -14256.75
$ws.Range("H97").Value = 2300
$ws.Range("J97").Value = 2300
$ws.Range("L97").Value = 6900
$ws.Range("N97").Value = -7892
$ws.Range("H112").Value = 3229.5
$ws.Range("I112").Value = 1331
$ws.Range("J112").Value = 3413.2258
$ws.Range("K112").Value = 3993
$ws.Range("L112").Value = 10239.6774
$ws.Range("M112").Value = -2885
$ws.Range("N112").Value = -12455.6774
$ws.Range("H113").Value = 1995
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1995
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1995
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -8503
$ws.Range("H137").Value = 2194.5
$ws.Range("I137").Value = 1625.25
$ws.Range("J137").Value = 2574
$ws.Range("K137").Value = 4875.75
$ws.Range("L137").Value = 7722
$ws.Range("M137").Value = -2325.75
$ws.Range("N137").Value = -12822
$ws.Range("H138").Value = 8406.341
$ws.Range("I138").Value = 4192
$ws.Range("J138").Value = 9270.82
$ws.Range("K138").Value = 12576
$ws.Range("L138").Value = 27812.46
$ws.Range("M138").Value = -7436
$ws.Range("N138").Value = -38092.46

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2757.4055
$ws.Range("J2").Value = 2871.375
$ws.Range("L2").Value = 2871.375
$ws.Range("N2").Value = -3097.375
$ws.Range("H45").Value = 3333.9644
$ws.Range("I45").Value = 2866.0625
$ws.Range("K45").Value = 2866.0625
$ws.Range("M45").Value = -2489.0625
$ws.Range("H74").Value = 1985.5555
$ws.Range("I74").Value = 1464.7142
$ws.Range("K74").Value = 1464.7142
$ws.Range("M74").Value = -590.7141999999999
$ws.Range("H77").Value = 1985.5555
$ws.Range("I77").Value = 1464.7142
$ws.Range("K77").Value = 7323.571
$ws.Range("M77").Value = -2955.571
$ws.Range("H116").Value = 2757.4055
$ws.Range("J116").Value = 2871.375
$ws.Range("L116").Value = 2871.375
$ws.Range("N116").Value = -7459.375
$ws.Range("H132").Value = 6044
$ws.Range("I132").Value = 6063.75
$ws.Range("J132").Value = 5991.3335
$ws.Range("K132").Value = 18191.25
$ws.Range("L132").Value = 17974.0005
$ws.Range("M132").Value = -15661.25
$ws.Range("N132").Value = -23034.0005
$ws.Range("H133").Value = 79556.60000000001
$ws.Range("J133").Value = 79556.60000000001
$ws.Range("L133").Value = 79556.60000000001
$ws.Range("N133").Value = -84616.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2757.4055
$ws.Range("J3").Value = 2871.375
$ws.Range("L3").Value = 2871.375
$ws.Range("N3").Value = -3099.375
$ws.Range("H134").Value = 998.5
$ws.Range("I134").Value = 998.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 2995.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -460.5
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1066.6666
$ws.Range("I22").Value = 474.5
$ws.Range("J22").Value = 2251
$ws.Range("K22").Value = 474.5
$ws.Range("L22").Value = 2251
$ws.Range("M22").Value = -124.5
$ws.Range("N22").Value = -2951
$ws.Range("H31").Value = 5042.189
$ws.Range("I31").Value = 4895.619
$ws.Range("J31").Value = 5234.5625
$ws.Range("K31").Value = 4895.619
$ws.Range("L31").Value = 5234.5625
$ws.Range("M31").Value = -4600.619
$ws.Range("N31").Value = -5824.5625
$ws.Range("H34").Value = 5042.189
$ws.Range("I34").Value = 4895.619
$ws.Range("J34").Value = 5234.5625
$ws.Range("K34").Value = 4895.619
$ws.Range("L34").Value = 5234.5625
$ws.Range("M34").Value = -4693.619
$ws.Range("N34").Value = -5638.5625
$ws.Range("H138").Value = 87495.5
$ws.Range("J138").Value = 87495.5
$ws.Range("L138").Value = 87495.5
$ws.Range("N138").Value = -97775.5
$ws.Range("H139").Value = 83899
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 83899
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 83899
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -94179

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1294.7142
$ws.Range("I5").Value = 1242.75
$ws.Range("J5").Value = 1364
$ws.Range("K5").Value = 3728.25
$ws.Range("L5").Value = 4092
$ws.Range("M5").Value = -3616.25
$ws.Range("N5").Value = -4316
$ws.Range("H56").Value = 9135.544
$ws.Range("I56").Value = 9135.544
$ws.Range("K56").Value = 9135.544
$ws.Range("M56").Value = -8605.544
$ws.Range("H68").Value = 3361.1143
$ws.Range("I68").Value = 2810
$ws.Range("J68").Value = 3432.2258
$ws.Range("K68").Value = 8430
$ws.Range("L68").Value = 10296.6774
$ws.Range("M68").Value = -7619
$ws.Range("N68").Value = -11918.6774
$ws.Range("H71").Value = 3361.1143
$ws.Range("I71").Value = 2810
$ws.Range("J71").Value = 3432.2258
$ws.Range("K71").Value = 25290
$ws.Range("L71").Value = 30890.0322
$ws.Range("M71").Value = -21234
$ws.Range("N71").Value = -39002.0322
$ws.Range("H112").Value = 6363
$ws.Range("I112").Value = 7727
$ws.Range("J112").Value = 4999
$ws.Range("K112").Value = 23181
$ws.Range("L112").Value = 14997
$ws.Range("M112").Value = -22073
$ws.Range("N112").Value = -17213
$ws.Range("H129").Value = 2484.7
$ws.Range("I129").Value = 566.1429000000001
$ws.Range("J129").Value = 6961.3335
$ws.Range("K129").Value = 1698.4287
$ws.Range("L129").Value = 20884.0005
$ws.Range("M129").Value = 3301.5713
$ws.Range("N129").Value = -30884.0005
$ws.Range("H130").Value = 6566
$ws.Range("I130").Value = 3599
$ws.Range("J130").Value = 12500
$ws.Range("K130").Value = 10797
$ws.Range("L130").Value = 37500
$ws.Range("M130").Value = -5777
$ws.Range("N130").Value = -47540
$ws.Range("H135").Value = 1294.7142
$ws.Range("I135").Value = 1242.75
$ws.Range("J135").Value = 1364
$ws.Range("K135").Value = 11184.75
$ws.Range("L135").Value = 12276
$ws.Range("M135").Value = -8649.75
$ws.Range("N135").Value = -17346
$ws.Range("H137").Value = 4721.5
$ws.Range("I137").Value = 4100.2
$ws.Range("J137").Value = 5757
$ws.Range("K137").Value = 12300.6
$ws.Range("L137").Value = 17271
$ws.Range("M137").Value = -7200.599999999999
$ws.Range("N137").Value = -27471

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 18995
$ws.Range("J98").Value = 18995
$ws.Range("L98").Value = 18995
$ws.Range("N98").Value = -24985
$ws.Range("H132").Value = 3379.1538
$ws.Range("I132").Value = 3459.3872
$ws.Range("J132").Value = 3068.25
$ws.Range("K132").Value = 10378.1616
$ws.Range("L132").Value = 9204.75
$ws.Range("M132").Value = -7848.161599999999
$ws.Range("N132").Value = -14264.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4558.769
$ws.Range("I136").Value = 4364
$ws.Range("J136").Value = 4997
$ws.Range("K136").Value = 13092
$ws.Range("L136").Value = 14991
$ws.Range("M136").Value = -10542
$ws.Range("N136").Value = -20091

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2886.0571
$ws.Range("I132").Value = 2687.75
$ws.Range("J132").Value = 5001.3335
$ws.Range("K132").Value = 8063.25
$ws.Range("L132").Value = 15004.0005
$ws.Range("M132").Value = -5533.25
$ws.Range("N132").Value = -20064.0005

Write-Output "All updates applied."